$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spreadsheet headers (GPLIM-2588):
#  A1: "Sample ID" -> "Specimen_Number"
#  F1: "T/N"       -> "SAMPLE_TYPE"
$ws.Range("A1").Value = "Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

# Move the active selection to F2 (matches the saved sheet view state)
$ws.Range("F2").Select() | Out-Null
